$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without Excel
# re-interpreting number-looking strings (e.g. "212.02") as
# real numbers. Forcing a text NumberFormat first, then restoring
# the default 'Normal' style after the write, keeps the visible
# formatting/style identical to the original cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value2 = '26.179.67'
$ws.Range("E2").Value2 = '  +0.18%  '

$ws.Range("D3").Value2 = '1.602.47'
$ws.Range("E3").Value2 = '  -0.37%  '

$ws.Range("E4").Value2 = '  +0.18%  '

Set-TextValue ($ws.Range("D5")) '212.02'
$ws.Range("E5").Value2 = '  -0.33%  '

$ws.Range("E7").Value2 = '  +0.10%  '

$ws.Range("E8").Value2 = '  -0.76%  '

$ws.Range("E9").Value2 = '  -0.75%  '

Set-TextValue ($ws.Range("D10")) '18.13'
$ws.Range("E10").Value2 = '  -1.51%  '

$ws.Range("E11").Value2 = '  +2.16%  '

$ws.Range("D12").Value2 = '1.823.21'
$ws.Range("E12").Value2 = '  -0.40%  '

$ws.Range("D13").Value2 = '1.601.76'
$ws.Range("E13").Value2 = '  -0.08%  '

Set-TextValue ($ws.Range("D15")) '0.517'
$ws.Range("E15").Value2 = '  +1.00%  '

$ws.Range("D16").Value2 = '26.183.89'
$ws.Range("E16").Value2 = '  +0.16%  '

Set-TextValue ($ws.Range("D17")) '61.01'
$ws.Range("E17").Value2 = '  +0.24%  '

$ws.Range("E18").Value2 = '  -0.61%  '

$ws.Range("E19").Value2 = '  +0.28%  '

Set-TextValue ($ws.Range("D20")) '204.28'
$ws.Range("E20").Value2 = '  +2.84%  '

$ws.Range("E21").Value2 = '  -0.06%  '

$ws.Range("E22").Value2 = '  -2.29%  '

$ws.Range("E23").Value2 = '  +0.34%  '

Set-TextValue ($ws.Range("D24")) '1.93'
$ws.Range("E24").Value2 = '  +12.17%  '

$ws.Range("E25").Value2 = '  +1.09%  '

$ws.Range("E26").Value2 = '  +0.17%  '

$ws.Range("E27").Value2 = '  -7.64%  '

$ws.Range("E28").Value2 = '  -0.07%  '

$ws.Range("E29").Value2 = '  +0.13%  '

Set-TextValue ($ws.Range("D30")) '0.0487'
$ws.Range("E30").Value2 = '  +2.26%  '

Set-TextValue ($ws.Range("D31")) '1.16'
$ws.Range("E31").Value2 = '  -0.75%  '

$ws.Range("E32").Value2 = '  -0.24%  '

$ws.Range("E33").Value2 = '  -4.30%  '

$ws.Range("E34").Value2 = '  -2.16%  '

$ws.Range("E35").Value2 = '  +0.16%  '

$ws.Range("D36").Value2 = '1.144.08'
$ws.Range("E36").Value2 = '  +3.44%  '

$ws.Range("E37").Value2 = '  +6.51%  '

$ws.Range("E38").Value2 = '  +0.32%  '

$ws.Range("E39").Value2 = '  -1.80%  '

$ws.Range("E40").Value2 = '  -0.22%  '

$ws.Range("E41").Value2 = '  -2.55%  '

$ws.Range("E42").Value2 = '  -2.51%  '

$ws.Range("E43").Value2 = '  +0.30%  '

$ws.Range("D44").Value2 = '1.737.86'
$ws.Range("E44").Value2 = '  -0.29%  '

Set-TextValue ($ws.Range("D45")) '92.14'
$ws.Range("E45").Value2 = '  -1.09%  '

$ws.Range("E46").Value2 = '  -2.88%  '

$ws.Range("E47").Value2 = '  +0.27%  '

$ws.Range("E49").Value2 = '  -0.48%  '

$ws.Range("E50").Value2 = '  +0.40%  '

$ws.Range("D51").Value2 = '0.0₇0949'
$ws.Range("E51").Value2 = '  -11.68%  '
